# Update New Orleans hotel/review workbook:
#  1. Insert a new "State" column into hotel_info (between Hotel_Name and City)
#     and populate it with "Louisiana" for the existing data row.
#  2. Reorder the worksheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a new column C (State) - shifts City/Zip/etc. one column to the right.
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# Move hotel_info so that it sits after review_info (review_info becomes first tab).
$hotelSheet.Move($null, $reviewSheet)
